$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Fill in the missing "passed" status for the existing last row
# ("Admins can view a page to view newly submitted apps")
$ws.Range("B18").Value = "passed"

# Add the two new test case rows
$ws.Range("A19").Value = "All users can see all approved apps"
$ws.Range("B19").Value = "passed"

$ws.Range("A20").Value = "Moderators and admins can delete comments"
$ws.Range("B20").Value = "passed"

# Update the view: scroll so row 4 is at the top, and move the active
# selection down past the newly added rows.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()
